# edit.ps1 - applies the "zusätzliche Dokumente + Text Korrektur" commit.
#
# 1) Swap the order of slides 5 ("Diagramme") and 6 ("Sprache") so that
#    "Sprache" now comes before "Diagramme".
# 2) Update the cached "last saved" date field (05.07.19 -> 06.07.2019)
#    on the slide master and every slide layout.
# 3) Small text corrections on the "Zielgruppe" slide and the "Konzept" slide.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Reorder slides: move the slide titled "Sprache" (currently #6) to
#    position #5, pushing "Diagramme" (currently #5) down to #6.
# ---------------------------------------------------------------------
function Get-SlideTitle($slide) {
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            return $shape.TextFrame.TextRange.Text
        }
    }
    return ""
}

$spracheIndex = -1
$diagrammeIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $title = Get-SlideTitle $p.Slides.Item($i)
    if ($title -eq "Sprache") { $spracheIndex = $i }
    if ($title -eq "Diagramme") { $diagrammeIndex = $i }
}

if ($spracheIndex -gt 0 -and $diagrammeIndex -gt 0 -and $spracheIndex -gt $diagrammeIndex) {
    $p.Slides.Item($spracheIndex).MoveTo($diagrammeIndex)
}

# ---------------------------------------------------------------------
# 2) Fix the cached date field text on the master + every layout.
# ---------------------------------------------------------------------
function Set-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -eq "05.07.19") {
                $shape.TextFrame.TextRange.Text = "06.07.2019"
            }
        }
    }
}

$master = $p.SlideMaster
Set-DateText $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Set-DateText $master.CustomLayouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 3) Text corrections.
# ---------------------------------------------------------------------

# Paragraph ranges include a trailing CR/LF paragraph-mark character for
# every paragraph except the very last one in the text frame, so trim
# that off before comparing.
function Trim-ParaText($text) {
    return $text.TrimEnd("`r", "`n")
}

# Slide "Zielgruppe": "Allgemeinheit" -> "Allgemeinbevölkerung" and
# expand the "Fokus: ..." bullet with "Schulabsolventen,".
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ((Get-SlideTitle $slide) -eq "Zielgruppe") {
        foreach ($shape in $slide.Shapes) {
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $tr = $shape.TextFrame.TextRange
                $paras = $tr.Paragraphs()
                for ($pi = 1; $pi -le $paras.Count; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    $ptext = Trim-ParaText $para.Text
                    if ($ptext -eq "Allgemeinheit") {
                        $para.Text = "Allgemeinbevölkerung"
                    } elseif ($ptext -eq "Fokus: Studierende und Jobeinsteiger (18-28 J.)") {
                        $para.Text = "Fokus: Schulabsolventen, Studierende und Jobeinsteiger (18-28 J.)"
                    }
                }
            }
        }
    }
}

# Slide "Konzept": "Vollzeit/Teilzeit" -> "Beschäftigungsverhältnis"
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ((Get-SlideTitle $slide) -eq "Konzept") {
        foreach ($shape in $slide.Shapes) {
            if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
                $tr = $shape.TextFrame.TextRange
                $paras = $tr.Paragraphs()
                for ($pi = 1; $pi -le $paras.Count; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    $ptext = Trim-ParaText $para.Text
                    if ($ptext -eq "Vollzeit/Teilzeit") {
                        $para.Text = "Beschäftigungsverhältnis"
                    }
                }
            }
        }
    }
}

Write-Output "done"
